# Revert "Powerpoint writer: consolidate text run nodes."
#
# Split the trailing-space-bearing word runs back into a word run and a
# separate single-space run, on the three slide titles that were
# consolidated.

$p = $ppt.ActivePresentation

# Slide 1 title: "Header with inline code"
#   "Header " + "with " + "inline code"  ->  "Header" + " " + "with" + " " + "inline code"
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Characters(1, 6).Text = "Header"
$tr1.Characters(8, 4).Text = "with"

# Slide 2 title: "Syntax highlighting"
#   "Syntax " + "highlighting"  ->  "Syntax" + " " + "highlighting"
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$tr2.Characters(1, 6).Text = "Syntax"

# Slide 3 title: "Two column slide"
#   "Two " + "column " + "slide"  ->  "Two" + " " + "column" + " " + "slide"
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Characters(1, 3).Text = "Two"
$tr3.Characters(5, 6).Text = "column"
